$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 19 (pushes old rows 19/20 down to 20/21) ---
$ws.Rows("19").Insert()

# --- Fix the border on the new row so it reuses the existing "item row" style family ---
$rowRange = $ws.Range("A19:Q19")
$rowRange.Borders.Item(9).Color = 13882323
$rowRange.Borders.Item(9).LineStyle = 1

# --- Re-create the merges for the new row 19 (mirrors rows 7-18 layout) ---
$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# --- Row heights: new item row matches others (25.5); totals row becomes 24.75 ---
$ws.Rows("19").RowHeight = 25.5
$ws.Rows("20").RowHeight = 24.75

# --- Populate the new item row (#13: XITHRONE 500MG 5 F.C.TAB.) ---
$ws.Range("A19").Value = 13
$ws.Range("C19").Value = "XITHRONE 500MG 5 F.C.TAB."
$ws.Range("H19").Value = "1:0"

$origFmtL = $ws.Range("L19").NumberFormat
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = "1"
$ws.Range("L19").NumberFormat = $origFmtL

$ws.Range("N19").Value = "86.00"

$origFmtP = $ws.Range("P19").NumberFormat
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = "86.0000"
$ws.Range("P19").NumberFormat = $origFmtP

$ws.Range("Q19").Value = "1:0"

# --- Update the totals cell (was 716.8, now 802.8) ---
$ws.Range("P20").Value = 802.79999999999995

# --- Update the footer timestamp (10:31 AM -> 10:32 AM) ---
$ws.Range("A21").Value = "Wednesday, 18 June, 2025 10:32 AM"
